# Updates the Coin/Link/Price/Volume(1h) table on Sheet1 (cryptos.xlsx)
# to the latest scrape, per GitHub Actions commit
# "Updated cryptos list on Thu Apr  6 02:22:32 UTC 2023 with GitHub Actions".
#
# The ranking shifted by one row starting at row 26 (LEO dropped out of the
# list and every following coin moved up a row, with EOS newly appended at
# the bottom), so B/C (name/link) are rewritten alongside D/E for rows 26-51.
#
# Column D prices are stored as plain text (not numbers) in the workbook -
# some of them look numeric (e.g. '1.000', '0.2210') and Excel would silently
# coerce them to numbers (dropping the meaningful trailing zeros) unless they
# are entered with a leading apostrophe, exactly like a user typing them in
# manually would do to force text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.081.80'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').Value = '1.898.43'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''312.89'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '''0.5027'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = '''0.3896'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('D9').Value = '''0.09159'
$ws.Range('E9').Value = '  -5.20%  '
$ws.Range('D10').Value = '''1.131'
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('D11').Value = '''41.75'
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '''6.380'
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').Value = '''20.82'
$ws.Range('D14').Value = '1.889.26'
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('D15').Value = '''7.286'
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '''92.45'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '''0.00001108'
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').Value = '''0.06654'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '''17.91'
$ws.Range('E20').Value = '  -0.61%  '
$ws.Range('D21').Value = '''0.9999'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = '''6.209'
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('D23').Value = '28.138.32'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('D24').Value = '''11.42'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  +1.23%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '''2.576'
$ws.Range('E26').Value = '  -6.89%  '
$ws.Range('D27').Value = '2.109.90'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''20.89'
$ws.Range('E28').Value = '  -2.96%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = '''157.96'
$ws.Range('E29').Value = '  -1.09%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''126.24'
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.090'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.1060'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''5.616'
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.616'
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').Value = '''9.589'
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.06619'
$ws.Range('E36').Value = '  -2.94%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02408'
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '''0.2210'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''1.225'
$ws.Range('E39').Value = '  -3.07%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''1.267'
$ws.Range('E40').Value = '  +5.64%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.6485'
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').Value = '''4.974'
$ws.Range('E42').Value = '  -2.36%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''11.42'
$ws.Range('E43').Value = '  -2.02%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '''0.9999'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.6095'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''13.33'
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '''1.297'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '''3.688'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''2.000'
$ws.Range('E49').Value = '  -2.14%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''121.99'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '''1.183'
$ws.Range('E51').Value = '  -2.50%  '
